$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update tube radius values (AArm_tube_Rint / AArm_tube_Rext)
$ws.Range("B6").Value = 8
$ws.Range("B7").Value = 9.75

# Update tube radius values (Suspension_Rod_Rint / Suspension_Rod_Rext)
$ws.Range("B20").Value = 8
$ws.Range("B21").Value = 9.75

# Update selected cell on sheet
$ws.Range("E10").Select()
